$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New timeline rows (row, date-serial, country, item, quantity).
# Quantity "Undisclosed" (row 55) is a text value, not a number.
$data = @(
    ,@(43, 45203, "Switzerland", "GCS-200 de-mining vehicle", 1)
    ,@(44, 45201, "Germany", "Bandvagn 206 ATV", 14)
    ,@(45, 45201, "Germany", "Beaver bridge-laying tank", 1)
    ,@(46, 45201, "Germany", "Wisent mine-clearing tank", 2)
    ,@(47, 45201, "Germany", "SatCom terminals", 99)
    ,@(48, 45201, "Germany", "VECTOR drone spare parts package", 1)
    ,@(49, 45201, "Germany", "Border patrol cars", 21)
    ,@(50, 45201, "Germany", "Safety glasses", 20000)
    ,@(51, 45201, "Germany", "Cryptophone (encrypted sat phone)", 239)
    ,@(52, 45201, "Germany", "HX81 tank transport tractor", 11)
    ,@(53, 45201, "Germany", "Semi-trailers", 12)
    ,@(54, 45201, "Germany", "40mm rounds", 32823)
    ,@(55, 45201, "Germany", "1202 infusion kits", "Undisclosed")
)

foreach ($row in $data) {
    $r = $row[0]
    $dateSerial = $row[1]
    $country = $row[2]
    $item = $row[3]
    $qty = $row[4]

    # Clone formatting (date number format, borders, etc.) from the row above
    # so the new row matches the rest of the table exactly.
    $srcRow = $r - 1
    $ws.Range("B" + $srcRow + ":F" + $srcRow).Copy($ws.Range("B" + $r + ":F" + $r))

    $ws.Range("B" + $r).Value = $dateSerial
    $ws.Range("C" + $r).Value = $country
    $ws.Range("D" + $r).Value = "Delivery"
    $ws.Range("E" + $r).Value = $item
    if ($qty -is [string]) {
        $ws.Range("F" + $r).Value = $qty
    } else {
        $ws.Range("F" + $r).Value = [double]$qty
    }
}

# Scroll/select as Excel would after entering the new rows.
[void]$ws.Range("A40").Select()
[void]$ws.Range("B56").Select()
